$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the existing table to its final B1:N3 footprint *before* inserting
# the new column, so the insert (which only touches column A, left of the
# table) doesn't disturb the table's column name bindings.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("B1:N3"))

# Insert a new column before column A - this shifts Team..Yellow Cards
# from A:M to B:N (and shifts the column width definitions with them).
$ws.Columns.Item(1).Insert()

# New "Season" column header + values for the two rows.
$ws.Range("A1").Value2 = "Season"
$ws.Range("A2").Value2 = "24/25"
$ws.Range("A3").Value2 = "23/24"

# New row 3 - Manchester City's 2023/24 season stats.
$ws.Range("B3").Value2 = "Manchester City"
$ws.Range("C3").Value2 = 1
$ws.Range("D3").Value2 = 91
$ws.Range("E3").Value2 = 38
$ws.Range("F3").Value2 = 28
$ws.Range("G3").Value2 = 3
$ws.Range("H3").Value2 = 7
$ws.Range("I3").Value2 = 96
$ws.Range("J3").Value2 = 34
$ws.Range("K3").Value2 = 62
$ws.Range("L3").Value2 = "Erling Haaland"
$ws.Range("M3").Value2 = 2
$ws.Range("N3").Value2 = 52

# Turn on a plain worksheet AutoFilter over the new Season column.
$ws.Range("A1:A3").AutoFilter()

# Excel records that filter as a hidden workbook-level defined name.
$name = $ws.Names.Add("_xlnm._FilterDatabase", "='Manchester City Stats'!`$A`$1:`$A`$3")
$name.Visible = $false

# Move the selection / view like the saved file shows.
$ws.Range("E10").Select()
